$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.897.90'
$ws.Range("E2").Value = '  -2.06%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.783.76'
$ws.Range("E3").Value = '  -1.90%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.008'
$ws.Range("E4").Value = '  +0.40%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '312.52'
$ws.Range("E5").Value = '  -0.82%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.008'
$ws.Range("E6").Value = '  +0.45%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5108'
$ws.Range("E7").Value = '  +0.93%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3776'
$ws.Range("E8").Value = '  -2.11%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07757'
$ws.Range("E9").Value = '  -9.34%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.43'
$ws.Range("E10").Value = '  -1.28%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.085'
$ws.Range("E11").Value = '  -1.87%  '
$ws.Range("B12").Value = 'BinanceUSD'
$ws.Range("C12").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.008'
$ws.Range("E12").Value = '  +0.42%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.173'
$ws.Range("E13").Value = '  -3.30%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.24'
$ws.Range("E14").Value = '  -3.65%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '1.780.80'
$ws.Range("E15").Value = '  -1.70%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.174'
$ws.Range("E16").Value = '  -4.30%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '91.78'
$ws.Range("E17").Value = '  -1.73%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001067'
$ws.Range("E18").Value = '  -6.89%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06500'
$ws.Range("E19").Value = '  -2.46%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.007'
$ws.Range("E20").Value = '  +0.44%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.06'
$ws.Range("E21").Value = '  -3.66%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.894'
$ws.Range("E22").Value = '  -2.84%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '27.944.48'
$ws.Range("E23").Value = '  -2.00%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.91'
$ws.Range("E24").Value = '  -4.27%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.260'
$ws.Range("E25").Value = '  -0.32%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '158.11'
$ws.Range("E26").Value = '  +0.07%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.24'
$ws.Range("E27").Value = '  -4.57%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.988.18'
$ws.Range("E28").Value = '  -1.76%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.328'
$ws.Range("E29").Value = '  -2.24%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '121.90'
$ws.Range("E30").Value = '  -3.26%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.1061'
$ws.Range("E31").Value = '  -1.31%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.030'
$ws.Range("E32").Value = '  -6.91%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.629'
$ws.Range("E33").Value = '  -1.35%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.440'
$ws.Range("E34").Value = '  -5.14%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.07025'
$ws.Range("E35").Value = '  -6.42%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02296'
$ws.Range("E36").Value = '  -2.50%  '
$ws.Range("B37").Value = 'Algorand'
$ws.Range("C37").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.2106'
$ws.Range("E37").Value = '  -5.27%  '
$ws.Range("B38").Value = 'FraxShare'
$ws.Range("C38").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.595'
$ws.Range("E38").Value = '  -1.44%  '
$ws.Range("B39").Value = 'Aptos'
$ws.Range("C39").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '11.48'
$ws.Range("E39").Value = '  +2.24%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.989'
$ws.Range("E40").Value = '  -4.25%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6093'
$ws.Range("E41").Value = '  -3.66%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.148'
$ws.Range("E42").Value = '  -3.33%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.336'
$ws.Range("E43").Value = '  -4.95%  '
$ws.Range("B44").Value = 'PancakeSwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.746'
$ws.Range("E44").Value = '  +0.13%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '12.99'
$ws.Range("E45").Value = '  -4.59%  '
$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5897'
$ws.Range("E46").Value = '  -0.41%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '125.66'
$ws.Range("E47").Value = '  +0.10%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.218'
$ws.Range("E48").Value = '  +1.68%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.888'
$ws.Range("E49").Value = '  -5.11%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06714'
$ws.Range("E50").Value = '  -4.02%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.038'
$ws.Range("E51").Value = '  -2.96%  '
